# Add data for 2022-12-02
# - Rename the sheet / update the "through" date from November 23 to November 24
# - Update the running-month header text to match
# - Bump / add the affected neighborhood x month cell counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) to reflect the new "through" date.
$ws.Name = "Through 2022-11-24"

# Update the column header label for the in-progress month.
$ws.Range("B1").Value = "November 2022 (through November 24)"

# Incremented existing counts.
$ws.Range("M2").Value = 6
$ws.Range("AI2").Value = 3
$ws.Range("AT2").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("X4").Value = 3
$ws.Range("X5").Value = 15
$ws.Range("BP6").Value = 2
$ws.Range("BP7").Value = 6
$ws.Range("X10").Value = 4
$ws.Range("X11").Value = 16
$ws.Range("M14").Value = 4
$ws.Range("BP17").Value = 3
$ws.Range("AI21").Value = 3
$ws.Range("X25").Value = 11
$ws.Range("BE25").Value = 11
$ws.Range("BE28").Value = 2
$ws.Range("X43").Value = 2
$ws.Range("X48").Value = 2
$ws.Range("M57").Value = 2
$ws.Range("M59").Value = 4
$ws.Range("X76").Value = 5
$ws.Range("X85").Value = 2

# Newly populated cells (previously empty).
$ws.Range("X31").Value = 1
$ws.Range("CA41").Value = 1
$ws.Range("M58").Value = 1
$ws.Range("AT64").Value = 1
$ws.Range("X65").Value = 1
$ws.Range("AI72").Value = 1
$ws.Range("M75").Value = 1
$ws.Range("X89").Value = 1
